$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  Name="NEYRA PEREIRA MONICA AGNES ALEXANDRA";        Total=125},
    @{Row=3;  Name="MOGOLLON MORON KARIN LISSET";                 Total=103},
    @{Row=4;  Name="SILVA ALVARADO EVELYN DE JESUS";               Total=96},
    @{Row=5;  Name="MORETO ESPINOZA CRISTIAN ESTEBAN";             Total=93},
    @{Row=6;  Name="TUANAMA PIZANGO ELIZABETH";                    Total=92},
    @{Row=7;  Name="TENE TRABUCCO GIAN PIERRE";                    Total=92},
    @{Row=8;  Name="SANCARRANCO SANCHEZ DE CRUZ GISSELA SHANI";    Total=90},
    @{Row=9;  Name="CHIROQUE YARLEQUE BETTY ELIZABETH";            Total=88},
    @{Row=10; Name="HERRERA JUAN MANUEL";                          Total=88},
    @{Row=11; Name="AGURTO TINEO CESIA JIMENA";                    Total=80}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Total
}
